$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row needs to be inserted before the current
# row 143, pushing the existing row 143..182 data down to 144..183.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with this week's data.
$ws.Cells.Item(143, 1).Value = 8
$ws.Cells.Item(143, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(143, 3).Value = "Coquimbo"
$ws.Cells.Item(143, 4).Value = 44463
$ws.Cells.Item(143, 5).Value = 4
$ws.Cells.Item(143, 6).Value = 100112032
$ws.Cells.Item(143, 7).Value = "Zapallo italiano"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 600
$ws.Cells.Item(143, 11).Value = 10000
$ws.Cells.Item(143, 12).Value = 11000
$ws.Cells.Item(143, 13).Value = 10500
$ws.Cells.Item(143, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(143, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(143, 16).Value = 210
$ws.Cells.Item(143, 17).Value = 50
$ws.Cells.Item(143, 18).Value = "Hortaliza"
